# test chat with file
# Populate the "内容" (content) column (C) for the indicator rows that were
# previously "无" (none) with the actual extracted values/details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "单次"
$ws.Range("C7").Value = "34:49"

$ws.Range("C8").Value = "设备信息指标：`n- 交易卡号：6222****4054`n- 交易账户：4301****1741`n- 对方账户：0019****0002`n- 交易国家或地区简称：CHN`n- 记账币种：人民币"

$ws.Range("C9").Value = "01-09 17:34:49"
$ws.Range("C10").Value = "它显示了交易发生的确切日期和时间"
$ws.Range("C11").Value = "794.97"
$ws.Range("C12").Value = "794.97"

$ws.Range("C13").Value = "开户时间：2025-01-09 17:34:49`n开户地点：网上银行`n交易卡号：6222****4054`n交易账户：4301****1741`n交易户名：朱晗`n记账金额：794.97`n记账币种：人民币`n对方账户：0019****0002`n对方户名：银联转账(云闪付)`n对方账户行别：上海银联电子支付服务有限公司"

$ws.Range("C17").Value = "网上银行用户"
$ws.Range("C18").Value = "收支详细信息"
$ws.Range("C19").Value = "6222****4054"
$ws.Range("C20").Value = "794.97"
$ws.Range("C21").Value = "794.97"

$ws.Range("C22").Value = "目标账户开户信息指标：`n- 开户时间：2025年1月9日17:34:49`n- 开户地点：网上银行`n- 交易卡号：6222****4054`n- 交易账户：4301****1741`n- 朱晗（交易者姓名）`n- 交易时间：2025年1月9日17:34:49`n- 业务摘要：无卡支付`n- 对方账户行别：上海银联电子支付服务有限公司`n- 对方账户：0019****0002`n- 银行代码：CHN`n- 交易金额：794.97`n- 记账币种：人民币`n- 记账金额：794.97"

$ws.Range("C23").Value = "未提供"
$ws.Range("C24").Value = "0019****0002"
$ws.Range("C25").Value = "开户日期"
$ws.Range("C26").Value = "朱晗的职业是网上银行交易"
$ws.Range("C27").Value = "不适用（因为没有提供与目标账户相关的受教育程度信息）"
$ws.Range("C28").Value = "对方面账户所有者注册的联系电话或电子邮件地址等"
$ws.Range("C29").Value = "否"
$ws.Range("C30").Value = "否"
